# Updated cryptos list on Sat Jan  6 13:59:48 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.004.62'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '2.238.38'
$ws.Range('E3').Value = '  -0.34%  '
$ws.Range('E4').Value = '  +0.26%  '
$ws.Range('D5').Value = '''305.58'
$ws.Range('E5').Value = '  -3.86%  '
$ws.Range('D6').Value = '''94.62'
$ws.Range('E6').Value = '  -6.18%  '
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('E9').Value = '  -3.49%  '
$ws.Range('D10').Value = '''34.73'
$ws.Range('E10').Value = '  -5.46%  '
$ws.Range('E11').Value = '  -2.13%  '
$ws.Range('E12').Value = '  -3.91%  '
$ws.Range('D14').Value = '2.579.80'
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '2.239.61'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '''0.822'
$ws.Range('E16').Value = '  -3.01%  '
$ws.Range('E17').Value = '  -4.62%  '
$ws.Range('D18').Value = '43.896.81'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('D19').Value = '0.0₃0962'
$ws.Range('E19').Value = '  -0.90%  '
$ws.Range('D20').Value = '''12.14'
$ws.Range('E20').Value = '  -7.79%  '
$ws.Range('E21').Value = '  -1.78%  '
$ws.Range('D22').Value = '''65.52'
$ws.Range('E22').Value = '  +0.32%  '
$ws.Range('D23').Value = '''236.32'
$ws.Range('E23').Value = '  +0.64%  '
$ws.Range('E24').Value = '  -5.33%  '
$ws.Range('E25').Value = '  -4.14%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('B27').Value = 'InjectiveProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D27').Value = '''38.02'
$ws.Range('E27').Value = '  +2.09%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''9.84'
$ws.Range('E28').Value = '  -5.18%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '''2.19'
$ws.Range('E29').Value = '  -0.57%  '
$ws.Range('E30').Value = '  -1.47%  '
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').Value = '''150.03'
$ws.Range('E32').Value = '  -5.73%  '
$ws.Range('D33').Value = '''0.0796'
$ws.Range('E33').Value = '  -5.60%  '
$ws.Range('D34').Value = '''2.60'
$ws.Range('E34').Value = '  -2.94%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = '''0.109'
$ws.Range('E36').Value = '  -2.16%  '
$ws.Range('E37').Value = '  +1.26%  '
$ws.Range('E38').Value = '  -8.70%  '
$ws.Range('D39').Value = '''15.02'
$ws.Range('E39').Value = '  -6.76%  '
$ws.Range('E40').Value = '  -7.36%  '
$ws.Range('E41').Value = '  -6.95%  '
$ws.Range('E42').Value = '  -5.85%  '
$ws.Range('D43').Value = '''1.01'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').Value = '1.740.76'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').Value = '''84.81'
$ws.Range('E45').Value = '  +4.06%  '
$ws.Range('E46').Value = '  -4.66%  '
$ws.Range('D47').Value = '''99.82'
$ws.Range('E47').Value = '  -2.11%  '
$ws.Range('E48').Value = '  -3.80%  '
$ws.Range('E49').Value = '  -1.50%  '
$ws.Range('D50').Value = '''68.69'
$ws.Range('E50').Value = '  -7.50%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = '''53.77'
$ws.Range('E51').Value = '  -6.14%  '
